$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Effort R 1.0")
$ws2 = $wb.Worksheets.Item("Effort R 0.9")

# --- Row 28: combine Effort [h] (B28=2) and Additional Effort [h] (C28=0.25) into B28=2.25 ---
$ws1.Range("B28").Value = 2.25
$ws1.Range("C28").ClearContents() | Out-Null

# --- Row 31: combine Effort [h] (B31=2.5) and Additional Effort [h] (C31=1.5) into B31=4 ---
$ws1.Range("B31").Value = 4
$ws1.Range("C31").ClearContents() | Out-Null

# --- New shared string used by the two new rows ---
# (writing the text directly lets Excel append it to the shared-string table)
$newTask = "Preparation of setup and release"

# --- New row 46 ---
$ws1.Range("A46").Value = 41486
$ws1.Range("B46").Value = 2
$ws1.Range("D46").Value = $newTask

# --- New row 47 ---
$ws1.Range("A47").Value = 41489
$ws1.Range("B47").Value = 1.5
$ws1.Range("D47").Value = $newTask

# --- Selections / active sheet ---
# Sheet "Effort R 0.9" selection changes to the full column A (not the active sheet)
[void]$ws2.Columns("A:A").Select()

# Sheet "Effort R 1.0" stays the tab-selected / active sheet, with C28 as the
# active cell
$ws1.Activate()
[void]$ws1.Range("C28").Select()
